# repull data, push all data, mean calculation
# Update the dSF (column F) values for rows whose source data was re-pulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -9
    10 = 0
    11 = -2
    12 = -2
    16 = -3
    17 = 8
    21 = -1
    26 = -2
    28 = -2
    31 = -2
    34 = -2
    40 = 0
    44 = -1
    50 = 0
    51 = 0
    52 = -1
    63 = 3
    65 = -1
    68 = 2
    73 = -1
    74 = 8
    79 = 0
    81 = 0
    82 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
